$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds free-text "price" strings (dotted thousands separators,
# trailing-zero decimals, etc.) that must stay literal text rather than be
# coerced into numbers by the normal Value-assignment parsing. Prefixing with
# an apostrophe forces text entry; resetting Style to "Normal" afterwards
# strips the quote-prefix cell style so no stray formatting is introduced.
function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "28.636.44"
$ws.Range("E2").Value = "  +0.88%  "
Set-TextValue "D3" "1.563.50"
$ws.Range("E3").Value = "  -0.14%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue "D5" "210.22"
$ws.Range("E5").Value = "  -0.32%  "
Set-TextValue "D6" "0.510"
$ws.Range("E6").Value = "  +4.08%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +4.89%  "
Set-TextValue "D9" "0.246"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("E10").Value = "  -0.40%  "
Set-TextValue "D11" "0.0898"
$ws.Range("E11").Value = "  +0.44%  "
Set-TextValue "D12" "1.787.77"
$ws.Range("E12").Value = "  -0.12%  "
Set-TextValue "D13" "1.562.40"
$ws.Range("E13").Value = "  -0.20%  "
Set-TextValue "D14" "28.679.43"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  -1.02%  "
Set-TextValue "D17" "61.57"
$ws.Range("E17").Value = "  +1.40%  "
Set-TextValue "D18" "227.68"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  +0.00%  "
Set-TextValue "D21" "0.999"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -0.55%  "
Set-TextValue "D23" "9.06"
$ws.Range("E23").Value = "  +1.21%  "
Set-TextValue "D24" "2.07"
$ws.Range("E24").Value = "  +0.88%  "
Set-TextValue "D25" "151.76"
$ws.Range("E25").Value = "  +0.88%  "
Set-TextValue "D26" "0.106"
$ws.Range("E26").Value = "  +2.86%  "
Set-TextValue "D27" "14.78"
$ws.Range("E27").Value = "  -0.79%  "
Set-TextValue "D28" "1.00"
Set-TextValue "D29" "6.25"
$ws.Range("E29").Value = "  -1.20%  "
Set-TextValue "D30" "0.0458"
$ws.Range("E30").Value = "  -3.92%  "
$ws.Range("E31").Value = "  -0.75%  "
Set-TextValue "D33" "1.403.69"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("E36").Value = "  -2.21%  "
Set-TextValue "D37" "2.67"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  -2.07%  "
Set-TextValue "D39" "0.0163"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.22%  "
Set-TextValue "D41" "0.517"
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("E44").Value = "  -2.45%  "
Set-TextValue "D45" "63.95"
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  -2.46%  "
Set-TextValue "D47" "1.699.64"
$ws.Range("E47").Value = "  -0.09%  "
Set-TextValue "D48" "0.854"
$ws.Range("E48").Value = "  -6.93%  "
Set-TextValue "D49" "84.71"
$ws.Range("E49").Value = "  -0.98%  "
Set-TextValue "D50" "42.60"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("E51").Value = "  -0.44%  "
